# Generate Report for Handback
# This script applies the "handback received / version mismatch" update for the
# 1eca0919-2244-4d90-ab37-e71fa61e620d file on both the zh-cn and de-de sheets,
# widens the "Error Detail" column, and adds a hyperlink on the newly filled in
# "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a846429600532cdbf65ca0f4d50c95ef4b3c6dce/e2e/1eca0919-2244-4d90-ab37-e71fa61e620d.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca249072854483dedada2480f888a408261af3f6/e2e/1eca0919-2244-4d90-ab37-e71fa61e620d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a846429600532cdbf65ca0f4d50c95ef4b3c6dce/e2e/1eca0919-2244-4d90-ab37-e71fa61e620d.md."

function Update-HandbackSheet($ws, $handbackFile, $handbackDatetime) {
    # Widen column P (Error Detail) to 40 characters.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # I6: Latest Target File -- same display text/target as the source file hyperlink (A6)
    $ws.Range("I6").Value = "1eca0919-2244-4d90-ab37-e71fa61e620d.md"
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestUrl, "", "", "1eca0919-2244-4d90-ab37-e71fa61e620d.md") | Out-Null
    $ws.Range("I6").Font.Name = "Calibri"
    $ws.Range("I6").Font.Size = 11
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # J6: Latest Handback File
    $ws.Range("J6").Value = $handbackFile

    # K6: Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDatetime

    # P6: Error Detail
    $ws.Range("P6").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackSheet $wsZhCn "1eca0919-2244-4d90-ab37-e71fa61e620d.1867400bc89572ced06a573f3454686bd4ef04ab.zh-cn.xlf" "2016-08-23 14:51:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackSheet $wsDeDe "1eca0919-2244-4d90-ab37-e71fa61e620d.1867400bc89572ced06a573f3454686bd4ef04ab.de-de.xlf" "2016-08-23 14:52:24"
